$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-427) holds the "Förändrad" date, stored as the
# serial date number 45205 (2023-10-06). Bump every one of these
# values by one day to 45206 (2023-10-07), matching the diff exactly.
$ws.Range("C2:C427").Value = 45206
